$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.480.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.897.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4912"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2927"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06695"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.892.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07329"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.174"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6663"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.448.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007844"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.138.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.340"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +13.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "192.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.113"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.482"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("E28").Value = "  +6.13%  "
$ws.Range("E29").Value = "  +4.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.330"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09156"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.058"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05173"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7387"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.101"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.715"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01808"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.676"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9243"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.041"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4387"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.911"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9934"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "68.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +19.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1364"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.585"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.031"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05844"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3917"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.75%  "
